# Auto-generated edit script replicating the XML diff against Behemoth_Profits workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1270.5714
$ws.Range("I11").Value = 1270.5714
$ws.Range("K11").Value = 1270.5714
$ws.Range("M11").Value = -1130.5714

$ws.Range("H12").Value = 317.86667
$ws.Range("I12").Value = 340.6154
$ws.Range("J12").Value = 170
$ws.Range("K12").Value = 340.6154
$ws.Range("L12").Value = 170
$ws.Range("M12").Value = -170.6154
$ws.Range("N12").Value = -510

$ws.Range("H28").Value = 963.0909
$ws.Range("I28").Value = 613.5714
$ws.Range("K28").Value = 613.5714
$ws.Range("M28").Value = -128.5714

$ws.Range("H51").Value = 7029.8
$ws.Range("I51").Value = 6718.4287
$ws.Range("J51").Value = 7302.25
$ws.Range("K51").Value = 6718.4287
$ws.Range("L51").Value = 7302.25
$ws.Range("M51").Value = -6234.4287
$ws.Range("N51").Value = -8270.25

$ws.Range("H82").Value = 8223.111000000001
$ws.Range("I82").Value = 495.125
$ws.Range("J82").Value = 70047
$ws.Range("K82").Value = 1485.375
$ws.Range("L82").Value = 210141
$ws.Range("M82").Value = -1079.375
$ws.Range("N82").Value = -210953

$ws.Range("H85").Value = 8223.111000000001
$ws.Range("I85").Value = 495.125
$ws.Range("J85").Value = 70047
$ws.Range("K85").Value = 1485.375
$ws.Range("L85").Value = 210141
$ws.Range("M85").Value = -81.375
$ws.Range("N85").Value = -212949

$ws.Range("H132").Value = 2148.913
$ws.Range("I132").Value = 2085.7368
$ws.Range("K132").Value = 6257.2104
$ws.Range("M132").Value = -3727.2104

$ws.Range("H138").Value = 2579.9539
$ws.Range("I138").Value = 993.8461
$ws.Range("J138").Value = 2976.4807
$ws.Range("K138").Value = 2981.5383
$ws.Range("L138").Value = 8929.4421
$ws.Range("M138").Value = 2158.4617
$ws.Range("N138").Value = -19209.4421

$ws.Range("H141").Value = 2314.3635
$ws.Range("I141").Value = 2281.7144
$ws.Range("K141").Value = 6845.1432
$ws.Range("M141").Value = -1665.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2057.7273
$ws.Range("I45").Value = 1550.7142
$ws.Range("K45").Value = 1550.7142
$ws.Range("M45").Value = -1173.7142

$ws.Range("H97").Value = 1134.6154
$ws.Range("I97").Value = 1152.24
$ws.Range("K97").Value = 1152.24
$ws.Range("M97").Value = -656.24

$ws.Range("H102").Value = 16979.475
$ws.Range("I102").Value = 18712.354
$ws.Range("J102").Value = 2250
$ws.Range("K102").Value = 18712.354
$ws.Range("L102").Value = 2250
$ws.Range("M102").Value = -17090.354
$ws.Range("N102").Value = -5494

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2437.75
$ws.Range("I20").Value = 2435.1304
$ws.Range("K20").Value = 2435.1304
$ws.Range("M20").Value = -2188.1304

$ws.Range("H86").Value = 3398.3
$ws.Range("I86").Value = 3227.7144
$ws.Range("J86").Value = 3796.3333
$ws.Range("K86").Value = 3227.7144
$ws.Range("L86").Value = 3796.3333
$ws.Range("M86").Value = -2104.7144
$ws.Range("N86").Value = -6042.3333

$ws.Range("H89").Value = 3398.3
$ws.Range("I89").Value = 3227.7144
$ws.Range("J89").Value = 3796.3333
$ws.Range("K89").Value = 16138.572
$ws.Range("L89").Value = 18981.6665
$ws.Range("M89").Value = -10522.572
$ws.Range("N89").Value = -30213.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1830.1
$ws.Range("I7").Value = 216.21428
$ws.Range("J7").Value = 5595.8335
$ws.Range("K7").Value = 216.21428
$ws.Range("L7").Value = 5595.8335
$ws.Range("M7").Value = -103.21428
$ws.Range("N7").Value = -5821.8335

$ws.Range("H22").Value = 353.4
$ws.Range("I22").Value = 353.4
$ws.Range("K22").Value = 353.4
$ws.Range("M22").Value = -3.399999999999977

$ws.Range("H99").Value = 3505.2666
$ws.Range("I99").Value = 4136.5
$ws.Range("J99").Value = 3084.4443
$ws.Range("K99").Value = 4136.5
$ws.Range("L99").Value = 3084.4443
$ws.Range("M99").Value = -2638.5
$ws.Range("N99").Value = -6080.4443

$ws.Range("H105").Value = 2347.0625
$ws.Range("I105").Value = 1965.7693
$ws.Range("J105").Value = 3999.3333
$ws.Range("K105").Value = 1965.7693
$ws.Range("L105").Value = 3999.3333
$ws.Range("M105").Value = -218.7692999999999
$ws.Range("N105").Value = -7493.3333

$ws.Range("H122").Value = 2133.2
$ws.Range("I122").Value = 2302.5715
$ws.Range("K122").Value = 6907.7145
$ws.Range("M122").Value = -4457.7145

$ws.Range("H126").Value = 3505.2666
$ws.Range("I126").Value = 4136.5
$ws.Range("J126").Value = 3084.4443
$ws.Range("K126").Value = 12409.5
$ws.Range("L126").Value = 9253.332900000001
$ws.Range("M126").Value = -9939.5
$ws.Range("N126").Value = -14193.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7181017.5
$ws.Range("J4").Value = 84091.414
$ws.Range("L4").Value = 252274.242
$ws.Range("N4").Value = -252498.242

$ws.Range("H38").Value = 31.416666
$ws.Range("J38").Value = 43.57143
$ws.Range("L38").Value = 130.71429
$ws.Range("N38").Value = -824.71429

$ws.Range("H47").Value = 1169.3334
$ws.Range("I47").Value = 1169.3334
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 3508.0002
$ws.Range("L47").Value = 0
$ws.Range("M47").ClearContents()

$ws.Range("H122").Value = 725.03125
$ws.Range("J122").Value = 1770.1
$ws.Range("L122").Value = 15930.9
$ws.Range("N122").Value = -20830.9

$ws.Range("H139").Value = 2772.724
$ws.Range("I139").Value = 3001.2856
$ws.Range("J139").Value = 2700
$ws.Range("K139").Value = 9003.856800000001
$ws.Range("L139").Value = 8100
$ws.Range("M139").Value = -3863.856800000001
$ws.Range("N139").Value = -18380

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 30004
$ws.Range("I5").Value = 30004
$ws.Range("K5").Value = 30004
$ws.Range("M5").Value = -29892

$ws.Range("H70").Value = 4919.2
$ws.Range("I70").Value = 4899
$ws.Range("K70").Value = 4899
$ws.Range("M70").Value = -4629

$ws.Range("H73").Value = 4919.2
$ws.Range("I73").Value = 4899
$ws.Range("K73").Value = 4899
$ws.Range("M73").Value = -3963

$ws.Range("H80").Value = 2481.7778
$ws.Range("I80").Value = 2767.8
$ws.Range("J80").Value = 2124.25
$ws.Range("K80").Value = 2767.8
$ws.Range("L80").Value = 2124.25
$ws.Range("M80").Value = -1769.8
$ws.Range("N80").Value = -4120.25

$ws.Range("H83").Value = 2481.7778
$ws.Range("I83").Value = 2767.8
$ws.Range("J83").Value = 2124.25
$ws.Range("K83").Value = 13839
$ws.Range("L83").Value = 10621.25
$ws.Range("M83").Value = -8847
$ws.Range("N83").Value = -20605.25

$ws.Range("H102").Value = 3912.5186
$ws.Range("I102").Value = 2649.4285
$ws.Range("K102").Value = 2649.4285
$ws.Range("M102").Value = -1027.4285

$ws.Range("H126").Value = 4181.091
$ws.Range("I126").Value = 3999
$ws.Range("J126").Value = 4499.75
$ws.Range("K126").Value = 11997
$ws.Range("L126").Value = 13499.25
$ws.Range("M126").Value = -9527
$ws.Range("N126").Value = -18439.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H80").Value = 30000
$ws.Range("I80").Value = 30000
$ws.Range("K80").Value = 30000
$ws.Range("M80").Value = -28877

$ws.Range("H82").Value = 1349.5
$ws.Range("J82").Value = 2000.3334
$ws.Range("L82").Value = 2000.3334
$ws.Range("N82").Value = -2722.3334

$ws.Range("H83").Value = 30000
$ws.Range("I83").Value = 30000
$ws.Range("K83").Value = 90000
$ws.Range("M83").Value = -84384

$ws.Range("H85").Value = 1349.5
$ws.Range("J85").Value = 2000.3334
$ws.Range("L85").Value = 2000.3334
$ws.Range("N85").Value = -4496.3334

$ws.Range("H100").Value = 2260.5386
$ws.Range("I100").Value = 2078.9
$ws.Range("K100").Value = 2078.9
$ws.Range("M100").Value = -1537.9

$ws.Range("H132").Value = 85106
$ws.Range("I132").Value = 45327.78
$ws.Range("J132").Value = 1000005
$ws.Range("K132").Value = 135983.34
$ws.Range("L132").Value = 3000015
$ws.Range("M132").Value = -133453.34
$ws.Range("N132").Value = -3005075

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 88000
$ws.Range("I2").Value = 94333.336
$ws.Range("K2").Value = 94333.336
$ws.Range("M2").Value = -94221.336
